$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-NumCell($sheet, $ref, $val, $fmt) {
    $c = $sheet.Range($ref)
    $c.NumberFormat = $fmt
    $c.Value = $val
}

function Set-TextCell($sheet, $ref, $text) {
    $c = $sheet.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.NumberFormat = "General"
}

# --- Header text edits (rich text runs within shared strings) ---
# A8: "Volume 30   Number  6" -> "...7"  (the trailing "6" run, char 21, len 1)
$a8 = $ws.Range("A8")
$a8Chars = $a8.Characters(21, 1)
$a8Chars.Text = "7"

# C9: "Report Covering the Week  2/6/2023  Through  2/12/2023"
#  -> "Report Covering the Week  2/13/2023  Through  2/19/2023"
# Apply the later (right-hand) substitution first so the earlier offset stays valid.
$c9 = $ws.Range("C9")
$c9Chars2 = $c9.Characters(46, 9)
$c9Chars2.Text = "2/19/2023"
$c9Chars1 = $c9.Characters(27, 8)
$c9Chars1.Text = "2/13/2023"

# --- Crime-statistics table updates (rows 14-29) ---
Set-TextCell $ws "C14" "0"
Set-NumCell $ws "D14" 1 '#,##0'
Set-NumCell $ws "E14" -100 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "G14" 2 '#,##0'
Set-NumCell $ws "H14" 0 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "J14" 2 '#,##0'
Set-NumCell $ws "K14" 0 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "C15" 2 '#,##0'
Set-NumCell $ws "D15" 1 '#,##0'
Set-NumCell $ws "E15" 100 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "F15" 2 '#,##0'
Set-NumCell $ws "G15" 4 '#,##0'
Set-NumCell $ws "H15" -50 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "I15" 5 '#,##0'
Set-NumCell $ws "J15" 10 '#,##0'
Set-NumCell $ws "K15" -50 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "L15" 400 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "M15" 25 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "N15" -58.333333333333 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "C16" 8 '#,##0'
Set-NumCell $ws "D16" 10 '#,##0'
Set-NumCell $ws "E16" -20 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "F16" 32 '#,##0'
Set-NumCell $ws "H16" 0 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "I16" 84 '#,##0'
Set-NumCell $ws "J16" 70 '#,##0'
Set-NumCell $ws "K16" 20 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "L16" 162.5 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "M16" 55.555555555555 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "N16" -69.892473118279 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "C17" 14 '#,##0'
Set-NumCell $ws "D17" 13 '#,##0'
Set-NumCell $ws "E17" 7.692307692307 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "F17" 51 '#,##0'
Set-NumCell $ws "G17" 42 '#,##0'
Set-NumCell $ws "H17" 21.428571428571 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "I17" 117 '#,##0'
Set-NumCell $ws "J17" 76 '#,##0'
Set-NumCell $ws "K17" 53.947368421052 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "L17" 42.682926829268 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "M17" 98.305084745762 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "N17" 10.377358490566 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "C18" 11 '#,##0'
Set-NumCell $ws "D18" 4 '#,##0'
Set-NumCell $ws "E18" 175 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "F18" 41 '#,##0'
Set-NumCell $ws "G18" 25 '#,##0'
Set-NumCell $ws "H18" 64 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "I18" 63 '#,##0'
Set-NumCell $ws "J18" 37 '#,##0'
Set-NumCell $ws "K18" 70.270270270270 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "L18" 152 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "M18" 10.526315789473 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "N18" -65.945945945946 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "C19" 18 '#,##0'
Set-NumCell $ws "D19" 19 '#,##0'
Set-NumCell $ws "E19" -5.263157894736 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "F19" 68 '#,##0'
Set-NumCell $ws "G19" 71 '#,##0'
Set-NumCell $ws "H19" -4.225352112676 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "I19" 121 '#,##0'
Set-NumCell $ws "J19" 130 '#,##0'
Set-NumCell $ws "K19" -6.923076923076 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "L19" 72.857142857142 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "M19" 95.161290322580 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "N19" 39.080459770114 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "C20" 9 '#,##0'
Set-NumCell $ws "D20" 10 '#,##0'
Set-NumCell $ws "E20" -10 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "F20" 45 '#,##0'
Set-NumCell $ws "G20" 51 '#,##0'
Set-NumCell $ws "H20" -11.764705882352 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "I20" 99 '#,##0'
Set-NumCell $ws "J20" 104 '#,##0'
Set-NumCell $ws "K20" -4.807692307692 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "L20" 175 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "M20" 253.571428571429 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "N20" -66.211604095563 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "C21" 62 '#,##0'
Set-NumCell $ws "D21" 58 '#,##0'
Set-NumCell $ws "E21" 6.896551724137 '#,##0.00;"-"#,##0.00'
Set-NumCell $ws "F21" 241 '#,##0'
Set-NumCell $ws "G21" 227 '#,##0'
Set-NumCell $ws "H21" 6.167400881057 '#,##0.00;"-"#,##0.00'
Set-NumCell $ws "I21" 491 '#,##0'
Set-NumCell $ws "J21" 429 '#,##0'
Set-NumCell $ws "K21" 14.452214452214 '#,##0.00;"-"#,##0.00'
Set-NumCell $ws "L21" 97.983870967741 '#,##0.00;"-"#,##0.00'
Set-NumCell $ws "M21" 85.984848484848 '#,##0.00;"-"#,##0.00'
Set-NumCell $ws "N21" -49.485596707818 '#,##0.00;"-"#,##0.00'
Set-NumCell $ws "C22" 1 '#,##0'
Set-NumCell $ws "F22" 1 '#,##0'
Set-NumCell $ws "H22" -50 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "I22" 1 '#,##0'
Set-NumCell $ws "K22" -50 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "L22" 0 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "M22" -50 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "C23" 7 '#,##0'
Set-NumCell $ws "D23" 10 '#,##0'
Set-NumCell $ws "E23" -30 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "F23" 20 '#,##0'
Set-NumCell $ws "G23" 31 '#,##0'
Set-NumCell $ws "H23" -35.483870967741 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "I23" 44 '#,##0'
Set-NumCell $ws "J23" 46 '#,##0'
Set-NumCell $ws "K23" -4.347826086956 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "L23" 46.666666666666 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "M23" 51.724137931034 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "C24" 43 '#,##0'
Set-NumCell $ws "D24" 45 '#,##0'
Set-NumCell $ws "E24" -4.444444444444 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "F24" 171 '#,##0'
Set-NumCell $ws "G24" 133 '#,##0'
Set-NumCell $ws "H24" 28.571428571428 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "I24" 284 '#,##0'
Set-NumCell $ws "J24" 213 '#,##0'
Set-NumCell $ws "K24" 33.333333333333 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "L24" 83.225806451612 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "M24" 56.906077348066 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "C25" 21 '#,##0'
Set-NumCell $ws "D25" 18 '#,##0'
Set-NumCell $ws "E25" 16.666666666666 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "F25" 79 '#,##0'
Set-NumCell $ws "G25" 80 '#,##0'
Set-NumCell $ws "H25" -1.25 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "I25" 149 '#,##0'
Set-NumCell $ws "J25" 135 '#,##0'
Set-NumCell $ws "K25" 10.370370370370 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "L25" 49 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "M25" -22.395833333333 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "C26" 2 '#,##0'
Set-NumCell $ws "E26" 100 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "F26" 2 '#,##0'
Set-NumCell $ws "G26" 5 '#,##0'
Set-NumCell $ws "H26" -60 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "I26" 7 '#,##0'
Set-NumCell $ws "J26" 12 '#,##0'
Set-NumCell $ws "K26" -41.666666666666 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "L26" 250 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "C27" 5 '#,##0'
Set-NumCell $ws "D27" 1 '#,##0'
Set-NumCell $ws "E27" 400 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "F27" 11 '#,##0'
Set-NumCell $ws "H27" 57.142857142857 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "I27" 15 '#,##0'
Set-NumCell $ws "J27" 9 '#,##0'
Set-NumCell $ws "K27" 66.666666666666 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "L27" 15.384615384615 '#,##0.0;"-"#,##0.0'
Set-TextCell $ws "C28" "0"
Set-TextCell $ws "D28" "0"
Set-TextCell $ws "E28" "***.*"
Set-NumCell $ws "F28" 2 '#,##0'
Set-NumCell $ws "H28" -50 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "L28" -57.142857142857 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "N28" -85 '#,##0.0;"-"#,##0.0'
Set-TextCell $ws "C29" "0"
Set-TextCell $ws "D29" "0"
Set-TextCell $ws "E29" "***.*"
Set-NumCell $ws "F29" 2 '#,##0'
Set-NumCell $ws "H29" -50 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "L29" -57.142857142857 '#,##0.0;"-"#,##0.0'
Set-NumCell $ws "N29" -83.333333333333 '#,##0.0;"-"#,##0.0'
